# Generate Report for Handoff
# Updates the localization-status report to reflect that a new handoff
# package was generated for b.md: status flips from "Handed back" to
# "Ready for handoff", the duplicate-content flag clears, new handoff
# file names / timestamps are recorded, and a staleness warning is added.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/c7f9934434264f51fb883a95351ab07c1353f5c6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/9b49fd2b6c048d23b4de5a3fd76cf8a963f84f1f/e2e/b.md."

# ---------------------------------------------------------------
# Overview sheet: update the b.md row (row 3) status + date columns
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-12 02:59:34"

# ---------------------------------------------------------------
# zh-cn sheet: widen the Error Detail column and update the b.md row
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-12 02:59:29"
$zhcn.Range("P3").Value = $errorDetail

# ---------------------------------------------------------------
# de-de sheet: widen the Error Detail column and update the b.md row
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-12 02:59:34"
$dede.Range("P3").Value = $errorDetail
